# This script re-applies a set of paired row-content swaps on the
# "Artfynd" worksheet (rows 4/5, 6/7, 8/9, 10/11, 12/13, 15/16).
# For every pair the "species record" data (columns A,B,D,E,F,G,H,Q,R
# and, when present, K,L,M,N,AF) is exchanged between the two rows,
# while all other columns (location name, county, dates, reporter, ...)
# stay untouched because they are identical for every row in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that always contain a value on every affected row.
$coreCols = @("A","B","D","E","F","G","H","Q","R")

function Swap-CoreColumns($ws, $rowA, $rowB, $cols) {
    $valsA = @{}
    $valsB = @{}
    foreach ($col in $cols) {
        $valsA[$col] = $ws.Range($col + $rowA).Value2
        $valsB[$col] = $ws.Range($col + $rowB).Value2
    }
    foreach ($col in $cols) {
        $ws.Range($col + $rowA).Value2 = $valsB[$col]
        $ws.Range($col + $rowB).Value2 = $valsA[$col]
    }
}

function Set-BlankCell($ws, $addr) {
    # Forces the cell to exist in the saved XML even though it has no
    # value (mirrors the workbook's original empty placeholder cells).
    $ws.Range($addr).NumberFormat = $ws.Range($addr).NumberFormat
}

function Clear-CellCompletely($ws, $addr) {
    $ws.Range($addr).ClearContents() | Out-Null
}

# rowA, rowB, and the optional columns (K,L,M,N,AF) that may exist only
# on one side. For each such column we give the value that row A should
# end up holding and the value row B should end up holding, using $null
# to mean "the cell must not exist at all".
$pairs = @(
    @{ A = 4;  B = 5;  Optional = @{ "AF" = @{ A = ""; B = $null } } },
    @{ A = 6;  B = 7;  Optional = @{ "AF" = @{ A = $null; B = "" } } },
    @{ A = 8;  B = 9;  Optional = @{
            "K" = @{ A = $null; B = "" }
            "L" = @{ A = $null; B = "" }
            "M" = @{ A = $null; B = "äldre spår" }
            "N" = @{ A = $null; B = "" }
        } },
    @{ A = 10; B = 11; Optional = @{} },
    @{ A = 12; B = 13; Optional = @{} },
    @{ A = 15; B = 16; Optional = @{} }
)

foreach ($pair in $pairs) {
    $rowA = $pair.A
    $rowB = $pair.B

    Swap-CoreColumns $ws $rowA $rowB $coreCols

    foreach ($col in $pair.Optional.Keys) {
        $target = $pair.Optional[$col]
        $addrA = $col + $rowA
        $addrB = $col + $rowB

        if ($null -eq $target.A) {
            Clear-CellCompletely $ws $addrA
        } elseif ($target.A -eq "") {
            Set-BlankCell $ws $addrA
        } else {
            $ws.Range($addrA).Value2 = $target.A
        }

        if ($null -eq $target.B) {
            Clear-CellCompletely $ws $addrB
        } elseif ($target.B -eq "") {
            Set-BlankCell $ws $addrB
        } else {
            $ws.Range($addrB).Value2 = $target.B
        }
    }
}
